# Auto update stock data
# Updates the "as of" date (column A) from 2025/11/07 to 2025/11/08 for every
# ticker block, plus the corresponding refreshed metric values that moved as
# part of the same data refresh (mostly column B "EBITDA", plus C74).
#
# All of these source cells are stored as literal text (t="inlineStr") rather
# than numbers/dates, so each write below first forces the Text number format
# on the target cell (so Excel doesn't auto-coerce "2025/11/08" into a date
# serial or "4.55" into a float), sets the literal value, and then restores
# the cell's original ("Normal") style so no stray number-format/style is
# left behind on cells that were previously unstyled (style index 0).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param(
        [string]$Address,
        [string]$Text
    )
    $rng = $ws.Range($Address)
    $rng.NumberFormat = "@"
    $rng.Value = $Text
    $rng.Style = "Normal"
}

# Row 2
Set-TextValue "A2" "2025/11/08"
Set-TextValue "B2" "4.55"

# Row 8
Set-TextValue "A8" "2025/11/08"
Set-TextValue "B8" "7.40"

# Row 14
Set-TextValue "A14" "2025/11/08"
Set-TextValue "B14" "2.74"

# Row 20
Set-TextValue "A20" "2025/11/08"
Set-TextValue "B20" "12.21"

# Row 26
Set-TextValue "A26" "2025/11/08"
Set-TextValue "B26" "9.88"

# Row 32
Set-TextValue "A32" "2025/11/08"
Set-TextValue "B32" "24.79"

# Row 38 (only the date changes; B38 stays 72.14)
Set-TextValue "A38" "2025/11/08"

# Row 44
Set-TextValue "A44" "2025/11/08"
Set-TextValue "B44" "11.21"

# Row 50
Set-TextValue "A50" "2025/11/08"
Set-TextValue "B50" "11.41"

# Row 56 (only the date changes; B56 stays blank)
Set-TextValue "A56" "2025/11/08"

# Row 62
Set-TextValue "A62" "2025/11/08"
Set-TextValue "B62" "11.61"

# Row 68
Set-TextValue "A68" "2025/11/08"
Set-TextValue "B68" "13.04"

# Row 74 (date + Debt/Equity Ratio column C; B74 stays 15.94)
Set-TextValue "A74" "2025/11/08"
Set-TextValue "C74" "1.76"
